$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A81").Value = "teste de jogo ps5"
$ws.Range("B81").Value = "Completo"
$ws.Range("C81").Value = "PS5"
$ws.Range("D81").Value = "Concluído"
